$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("getDataEntities")
$ws2 = $wb.Worksheets.Item("getDataGraphQL")

# Insert 3 new rows at the top of the data (rows 2-4), pushing existing
# test rows down. The inserted rows inherit blank default formatting, so
# copy the row-3-below-insertion-point's per-cell format onto them
# afterwards to match the sheet's established style (border + font).
$ws2.Rows("2:4").Insert()

# --- Row 2: JinZu-ApiEngine-Test-1 ---
$ws2.Range("A2").Value = 'JinZu-ApiEngine-Test-1'
$ws2.Range("B2").Value = 'good request, data retrieved'
$ws2.Range("C2").Value = '{Project(cond:"{status:{_eq:\"online\"},Lease_Group:{lease_type:{_eq:\"2\"}}}",order:"") {business_mgr business_unit charge_frequency city province district class_level classification_level credit_amount detail_address discount_ratio expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture no status name risk_mgr rent_type invert_Customer(cond:"",order:"") { actual_controller category cid city cname contact contact_detail ctype district enterprise_size group holding_type id is_connected_tx is_gov_fin_customer is_group_customer legal_person_id legal_person major_class middle_class office_address project province registered_address small_class} Restricted_By_Contract(cond:"",order:"") {accumulated_amount charge_frequency contract_amount customer grant_loan_frequency id lease_balance lease_end_time lease_num lease_start_time lease_unit leasing_principal make_loan_day overdue_amount overdue_days overdue_interest overdue_principal payment_method project} Refer_To_Lease_Group(cond:"",order:"") {asset_type count discount_ratio id lease_net_val lease_type lease_type_gb lease_type_yj nominal_cost project transfer_price unit_price}} }'
$ws2.Range("D2").Value = 200
$ws2.Range("E2").Value = 100000
$ws2.Range("F2").Value = 'Successfully'

# --- Row 3: JinZu-ApiEngine-Test-2 ---
$ws2.Range("A3").Value = 'JinZu-ApiEngine-Test-2'
$ws2.Range("B3").Value = 'good request, data retrieved'
$ws2.Range("C3").Value = '{Site(cond:"{id:{_eq:\"P000000666\"}}",order:"") { id location commissioning_date state power_station Has_Device_Inverter{ site pr production name type full_generation_hours} }}'
$ws2.Range("D3").Value = 200
$ws2.Range("E3").Value = 100000
$ws2.Range("F3").Value = 'Successfully'

# --- Row 4: JinZu-ApiEngine-Test-3 ---
$ws2.Range("A4").Value = 'JinZu-ApiEngine-Test-3'
$ws2.Range("B4").Value = 'good request, data retrieved'
$ws2.Range("C4").Value = '{Contract(cond:"{project:{_eq:33}}",order:"") { accumulated_amount charge_frequency contract_amount customer grant_loan_frequency id lease_balance lease_end_time lease_num lease_start_time lease_unit leasing_principal make_loan_day overdue_amount overdue_days overdue_interest overdue_principal payment_method project invert_Project(cond:"",order:"") { business_mgr business_unit charge_frequency city province district class_level classification_level credit_amount detail_address discount_ratio expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture no status name risk_mgr rent_type Refer_To_Lease_Group(cond:"",order:"") { asset_type count discount_ratio id lease_net_val lease_type lease_type_gb lease_type_yj nominal_cost project transfer_price unit_price Refer_To_Power_Station_Properties(cond:"",order:"") { ps_type structure avg_annual_eq_hours capacity } }  } invert_Customer(cond:"",order:"") { actual_controller category cid city cname contact contact_detail ctype district enterprise_size group holding_type id is_connected_tx is_gov_fin_customer is_group_customer legal_person_id legal_person major_class middle_class office_address project province registered_address small_class } } }
    '
$ws2.Range("D4").Value = 200
$ws2.Range("E4").Value = 100000
$ws2.Range("F4").Value = 'Successfully'

# Match formatting of the surrounding data rows (thin border, 10pt font)
$ws2.Range("A5:F5").Copy()
$ws2.Range("A2:F2").PasteSpecial(-4122)
$ws2.Range("A5:F5").Copy()
$ws2.Range("A3:F3").PasteSpecial(-4122)
$ws2.Range("A5:F5").Copy()
$ws2.Range("A4:F4").PasteSpecial(-4122)

# Row 4's query cell wraps text (matches other long-query cells elsewhere
# in the workbook that use the wrap-text style).
$ws3 = $wb.Worksheets.Item("queryJinzuByGraphQL")
$ws3.Range("B2").Copy()
$ws2.Range("C4").PasteSpecial(-4122)

$ws2.Rows("4:4").RowHeight = 13

# Re-set the cell values after the format paste (PasteSpecial(Formats)
# should not disturb them, but make sure nothing was clobbered).
$ws2.Range("A2").Value = 'JinZu-ApiEngine-Test-1'
$ws2.Range("C2").Value = '{Project(cond:"{status:{_eq:\"online\"},Lease_Group:{lease_type:{_eq:\"2\"}}}",order:"") {business_mgr business_unit charge_frequency city province district class_level classification_level credit_amount detail_address discount_ratio expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture no status name risk_mgr rent_type invert_Customer(cond:"",order:"") { actual_controller category cid city cname contact contact_detail ctype district enterprise_size group holding_type id is_connected_tx is_gov_fin_customer is_group_customer legal_person_id legal_person major_class middle_class office_address project province registered_address small_class} Restricted_By_Contract(cond:"",order:"") {accumulated_amount charge_frequency contract_amount customer grant_loan_frequency id lease_balance lease_end_time lease_num lease_start_time lease_unit leasing_principal make_loan_day overdue_amount overdue_days overdue_interest overdue_principal payment_method project} Refer_To_Lease_Group(cond:"",order:"") {asset_type count discount_ratio id lease_net_val lease_type lease_type_gb lease_type_yj nominal_cost project transfer_price unit_price}} }'
$ws2.Range("A3").Value = 'JinZu-ApiEngine-Test-2'
$ws2.Range("C3").Value = '{Site(cond:"{id:{_eq:\"P000000666\"}}",order:"") { id location commissioning_date state power_station Has_Device_Inverter{ site pr production name type full_generation_hours} }}'
$ws2.Range("A4").Value = 'JinZu-ApiEngine-Test-3'
$ws2.Range("C4").Value = '{Contract(cond:"{project:{_eq:33}}",order:"") { accumulated_amount charge_frequency contract_amount customer grant_loan_frequency id lease_balance lease_end_time lease_num lease_start_time lease_unit leasing_principal make_loan_day overdue_amount overdue_days overdue_interest overdue_principal payment_method project invert_Project(cond:"",order:"") { business_mgr business_unit charge_frequency city province district class_level classification_level credit_amount detail_address discount_ratio expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture no status name risk_mgr rent_type Refer_To_Lease_Group(cond:"",order:"") { asset_type count discount_ratio id lease_net_val lease_type lease_type_gb lease_type_yj nominal_cost project transfer_price unit_price Refer_To_Power_Station_Properties(cond:"",order:"") { ps_type structure avg_annual_eq_hours capacity } }  } invert_Customer(cond:"",order:"") { actual_controller category cid city cname contact contact_detail ctype district enterprise_size group holding_type id is_connected_tx is_gov_fin_customer is_group_customer legal_person_id legal_person major_class middle_class office_address project province registered_address small_class } } }
    '

# The previously-selected sheet (getDataEntities) loses tabSelected and
# getDataGraphQL becomes the active tab, with the selection resting on
# the newly added query cell.
$ws2.Activate()
$ws2.Range("C4").Select()
